$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 596.8125
$ws.Cells.Item(28, 9).Value = 615.8889
$ws.Cells.Item(28, 10).Value = 572.2857
$ws.Cells.Item(28, 11).Value = 615.8889
$ws.Cells.Item(28, 12).Value = 572.2857
$ws.Cells.Item(28, 13).Value = -130.8889
$ws.Cells.Item(28, 14).Value = -1542.2857

$ws.Cells.Item(113, 8).Value = 5641.25
$ws.Cells.Item(113, 9).Value = 4662.5
$ws.Cells.Item(113, 10).Value = 6620
$ws.Cells.Item(113, 11).Value = 4662.5
$ws.Cells.Item(113, 12).Value = 6620
$ws.Cells.Item(113, 13).Value = -1408.5
$ws.Cells.Item(113, 14).Value = -13128

$ws.Cells.Item(116, 8).Value = 1845.6666
$ws.Cells.Item(116, 9).Value = 1741
$ws.Cells.Item(116, 10).Value = 1873.2106
$ws.Cells.Item(116, 11).Value = 1741
$ws.Cells.Item(116, 12).Value = 1873.2106
$ws.Cells.Item(116, 13).Value = 1701
$ws.Cells.Item(116, 14).Value = -8757.2106

$ws.Cells.Item(137, 8).Value = 1853.4445
$ws.Cells.Item(137, 9).Value = 1307.125
$ws.Cells.Item(137, 10).Value = 2083.4736
$ws.Cells.Item(137, 11).Value = 3921.375
$ws.Cells.Item(137, 12).Value = 6250.4208
$ws.Cells.Item(137, 13).Value = -1371.375
$ws.Cells.Item(137, 14).Value = -11350.4208

$ws.Cells.Item(138, 8).Value = 3434.6726
$ws.Cells.Item(138, 9).Value = 1495.6522
$ws.Cells.Item(138, 10).Value = 4828.3438
$ws.Cells.Item(138, 11).Value = 4486.9566
$ws.Cells.Item(138, 12).Value = 14485.0314
$ws.Cells.Item(138, 13).Value = 653.0434000000005
$ws.Cells.Item(138, 14).Value = -24765.0314

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3952.3333
$ws.Cells.Item(61, 9).Value = 4154.4443
$ws.Cells.Item(61, 10).Value = 3649.1667
$ws.Cells.Item(61, 11).Value = 4154.4443
$ws.Cells.Item(61, 12).Value = 3649.1667
$ws.Cells.Item(61, 13).Value = -3942.4443
$ws.Cells.Item(61, 14).Value = -4073.1667

$ws.Cells.Item(74, 8).Value = 1430.9512
$ws.Cells.Item(74, 9).Value = 1204.75
$ws.Cells.Item(74, 10).Value = 1750.2941
$ws.Cells.Item(74, 11).Value = 1204.75
$ws.Cells.Item(74, 12).Value = 1750.2941
$ws.Cells.Item(74, 13).Value = -330.75
$ws.Cells.Item(74, 14).Value = -3498.2941

$ws.Cells.Item(77, 8).Value = 1430.9512
$ws.Cells.Item(77, 9).Value = 1204.75
$ws.Cells.Item(77, 10).Value = 1750.2941
$ws.Cells.Item(77, 11).Value = 6023.75
$ws.Cells.Item(77, 12).Value = 8751.470499999999
$ws.Cells.Item(77, 13).Value = -1655.75
$ws.Cells.Item(77, 14).Value = -17487.4705

$ws.Cells.Item(96, 8).Value = 47332.668
$ws.Cells.Item(96, 10).Value = 47332.668
$ws.Cells.Item(96, 12).Value = 47332.668
$ws.Cells.Item(96, 14).Value = -52824.668

$ws.Cells.Item(122, 8).Value = 1710959.2
$ws.Cells.Item(122, 9).Value = 2850141.5
$ws.Cells.Item(122, 10).Value = 2185.6667
$ws.Cells.Item(122, 11).Value = 8550424.5
$ws.Cells.Item(122, 12).Value = 6557.000100000001
$ws.Cells.Item(122, 13).Value = -8547974.5
$ws.Cells.Item(122, 14).Value = -11457.0001

$ws.Cells.Item(136, 8).Value = 3952.3333
$ws.Cells.Item(136, 9).Value = 4154.4443
$ws.Cells.Item(136, 10).Value = 3649.1667
$ws.Cells.Item(136, 11).Value = 12463.3329
$ws.Cells.Item(136, 12).Value = 10947.5001
$ws.Cells.Item(136, 13).Value = -9913.332900000001
$ws.Cells.Item(136, 14).Value = -16047.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(61, 8).Value = 0
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 14).Value = 0
$ws.Cells.Item(61, 12).ClearContents()

$ws.Cells.Item(86, 8).Value = 2439.6
$ws.Cells.Item(86, 9).Value = 1799.5
$ws.Cells.Item(86, 10).Value = 5000
$ws.Cells.Item(86, 11).Value = 1799.5
$ws.Cells.Item(86, 12).Value = 5000
$ws.Cells.Item(86, 13).Value = -676.5
$ws.Cells.Item(86, 14).Value = -7246

$ws.Cells.Item(89, 8).Value = 2439.6
$ws.Cells.Item(89, 9).Value = 1799.5
$ws.Cells.Item(89, 10).Value = 5000
$ws.Cells.Item(89, 11).Value = 8997.5
$ws.Cells.Item(89, 12).Value = 25000
$ws.Cells.Item(89, 13).Value = -3381.5
$ws.Cells.Item(89, 14).Value = -36232

$ws.Cells.Item(138, 8).Value = 57845.715
$ws.Cells.Item(138, 10).Value = 57845.715
$ws.Cells.Item(138, 12).Value = 57845.715
$ws.Cells.Item(138, 14).Value = -68125.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3758.75
$ws.Cells.Item(16, 9).Value = 1750
$ws.Cells.Item(16, 10).Value = 5767.5
$ws.Cells.Item(16, 11).Value = 1750
$ws.Cells.Item(16, 12).Value = 5767.5
$ws.Cells.Item(16, 13).Value = -1463
$ws.Cells.Item(16, 14).Value = -6341.5

$ws.Cells.Item(35, 8).Value = 6175
$ws.Cells.Item(35, 9).Value = 6175
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 6175
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 14).Value = -5881
$ws.Cells.Item(35, 13).ClearContents()

$ws.Cells.Item(58, 8).Value = 1471.7778
$ws.Cells.Item(58, 9).Value = 948.1667
$ws.Cells.Item(58, 10).Value = 2519
$ws.Cells.Item(58, 11).Value = 948.1667
$ws.Cells.Item(58, 12).Value = 2519
$ws.Cells.Item(58, 13).Value = -745.1667
$ws.Cells.Item(58, 14).Value = -2925

$ws.Cells.Item(112, 8).Value = 70702
$ws.Cells.Item(112, 10).Value = 70702
$ws.Cells.Item(112, 12).Value = 70702
$ws.Cells.Item(112, 14).Value = -73656

$ws.Cells.Item(113, 8).Value = 3758.75
$ws.Cells.Item(113, 9).Value = 1750
$ws.Cells.Item(113, 10).Value = 5767.5
$ws.Cells.Item(113, 11).Value = 1750
$ws.Cells.Item(113, 12).Value = 5767.5
$ws.Cells.Item(113, 13).Value = 420
$ws.Cells.Item(113, 14).Value = -10107.5

$ws.Cells.Item(136, 8).Value = 1471.7778
$ws.Cells.Item(136, 9).Value = 948.1667
$ws.Cells.Item(136, 10).Value = 2519
$ws.Cells.Item(136, 11).Value = 2844.5001
$ws.Cells.Item(136, 12).Value = 7557
$ws.Cells.Item(136, 13).Value = -294.5001000000002
$ws.Cells.Item(136, 14).Value = -12657

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 533333.3
$ws.Cells.Item(4, 9).Value = 533333.3
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 1599999.9
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 14).Value = -1599887.9
$ws.Cells.Item(4, 13).ClearContents()

$ws.Cells.Item(20, 8).Value = 5532.222
$ws.Cells.Item(20, 9).Value = 3700
$ws.Cells.Item(20, 10).Value = 6998
$ws.Cells.Item(20, 11).Value = 11100
$ws.Cells.Item(20, 12).Value = 20994
$ws.Cells.Item(20, 13).Value = -10873
$ws.Cells.Item(20, 14).Value = -21448

$ws.Cells.Item(107, 8).Value = 785.825
$ws.Cells.Item(107, 9).Value = 248.3125
$ws.Cells.Item(107, 11).Value = 744.9375
$ws.Cells.Item(107, 13).Value = 1175.0625

$ws.Cells.Item(131, 8).Value = 18645156
$ws.Cells.Item(131, 9).Value = 10000498
$ws.Cells.Item(131, 10).Value = 20409370
$ws.Cells.Item(131, 11).Value = 30001494
$ws.Cells.Item(131, 12).Value = 61228110
$ws.Cells.Item(131, 13).Value = -29996454
$ws.Cells.Item(131, 14).Value = -61238190

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 5349.3706
$ws.Cells.Item(126, 9).Value = 5617.32
$ws.Cells.Item(126, 10).Value = 2000
$ws.Cells.Item(126, 11).Value = 16851.96
$ws.Cells.Item(126, 12).Value = 6000
$ws.Cells.Item(126, 13).Value = -14381.96
$ws.Cells.Item(126, 14).Value = -10940

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(21, 8).Value = 52336.832
$ws.Cells.Item(21, 9).Value = 34000
$ws.Cells.Item(21, 11).Value = 34000
$ws.Cells.Item(21, 13).Value = -33826

$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 14).Value = 0
$ws.Cells.Item(45, 12).ClearContents()
$ws.Cells.Item(45, 13).ClearContents()

$ws.Cells.Item(68, 8).Value = 29413538
$ws.Cells.Item(68, 9).Value = 1629.6923
$ws.Cells.Item(68, 10).Value = 125002250
$ws.Cells.Item(68, 11).Value = 1629.6923
$ws.Cells.Item(68, 12).Value = 125002250
$ws.Cells.Item(68, 13).Value = -880.6922999999999
$ws.Cells.Item(68, 14).Value = -125003748

$ws.Cells.Item(71, 8).Value = 29413538
$ws.Cells.Item(71, 9).Value = 1629.6923
$ws.Cells.Item(71, 10).Value = 125002250
$ws.Cells.Item(71, 11).Value = 8148.461499999999
$ws.Cells.Item(71, 12).Value = 625011250
$ws.Cells.Item(71, 13).Value = -4404.461499999999
$ws.Cells.Item(71, 14).Value = -625018738

$ws.Cells.Item(132, 8).Value = 14448494
$ws.Cells.Item(132, 9).Value = 18844330
$ws.Cells.Item(132, 10).Value = 5037.2856
$ws.Cells.Item(132, 11).Value = 56532990
$ws.Cells.Item(132, 12).Value = 15111.8568
$ws.Cells.Item(132, 13).Value = -56530460
$ws.Cells.Item(132, 14).Value = -20171.8568

$ws.Cells.Item(139, 8).Value = 63290.668
$ws.Cells.Item(139, 10).Value = 63290.668
$ws.Cells.Item(139, 12).Value = 63290.668
$ws.Cells.Item(139, 14).Value = -73570.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(19, 8).Value = 63672.332
$ws.Cells.Item(19, 9).Value = 51005
$ws.Cells.Item(19, 11).Value = 51005
$ws.Cells.Item(19, 13).Value = -50831
